# Adds a new "2022-Q4" quarter to the 601808-中海油服 workbook:
#   - inserts a new "2022-Q4" worksheet (fund holdings detail) right after "总计"
#   - inserts a new summary row for 2022-Q4 at the top of the "总计" sheet's data
#
# Helper: force a cell's value to be stored as TEXT (so numeric-looking
# strings like fund codes "009394" or percentages "93.80" are not silently
# coerced into numbers / lose leading zeros / trailing zeros), then drop the
# "quote prefix" formatting that trick leaves behind so the cell keeps
# whatever style it already had (here: no explicit style, matching the
# source data).
function Set-TextCell($cell, $text) {
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$totalWs = $wb.Worksheets.Item(1)          # "总计"
$lastWs  = $wb.Worksheets.Item($wb.Worksheets.Count)   # original last tab ("2020-Q4")

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by duplicating the "2022-Q3" sheet
#    (index 2) right after "总计" -- this gives us the identical column
#    layout/styles (A1:H11, header + index-column formatting) for free.
# ---------------------------------------------------------------------
$q3Ws = $wb.Worksheets.Item(2)
$q3Ws.Copy($null, $totalWs)
$q4Ws = $wb.Worksheets.Item(2)
$q4Ws.Name = "2022-Q4"

# New fund-holdings data for 2022-Q4 (row order matches final column ranks).
$q4Rows = @(
    @{ B = "161810"; C = "银华内需精选混合（LOF）";         D = "23.02"; E = "93.80"; F = "6.70"; G = "1.5423"; H = 7 },
    @{ B = "009394"; C = "银华同力精选混合";                 D = "18.05"; E = "94.33"; F = "5.51"; G = "0.9946"; H = 9 },
    @{ B = "004641"; C = "万家量化睿选灵活配置混合A";        D = "9.45";  E = "89.25"; F = "1.36"; G = "0.1285"; H = 7 },
    @{ B = "011429"; C = "前海开源民裕进取混合";             D = "2.33";  E = "62.27"; F = "4.38"; G = "0.1021"; H = 3 },
    @{ B = "008480"; C = "永赢股息优选混合A";                D = "1.96";  E = "83.68"; F = "4.67"; G = "0.0915"; H = 6 },
    @{ B = "016556"; C = "万家量化睿选灵活配置混合C";        D = "4.91";  E = "89.25"; F = "1.36"; G = "0.0668"; H = 7 },
    @{ B = "080001"; C = "长盛成长价值混合A";                D = "2.43";  E = "63.89"; F = "2.68"; G = "0.0651"; H = 5 },
    @{ B = "011588"; C = "前海开源成份精选混合";             D = "0.77";  E = "62.53"; F = "4.57"; G = "0.0352"; H = 2 },
    @{ B = "008481"; C = "永赢股息优选混合C";                D = "0.24";  E = "83.68"; F = "4.67"; G = "0.0112"; H = 6 },
    @{ B = "012715"; C = "长盛成长价值混合C";                D = "0.10";  E = "63.89"; F = "2.68"; G = "0.0027"; H = 5 }
)

for ($i = 0; $i -lt $q4Rows.Count; $i++) {
    $r = $i + 2
    $row = $q4Rows[$i]
    Set-TextCell $q4Ws.Cells.Item($r, 2) $row.B   # fund code
    Set-TextCell $q4Ws.Cells.Item($r, 3) $row.C   # fund name
    Set-TextCell $q4Ws.Cells.Item($r, 4) $row.D   # fund size
    Set-TextCell $q4Ws.Cells.Item($r, 5) $row.E   # equity position %
    Set-TextCell $q4Ws.Cells.Item($r, 6) $row.F   # position ratio %
    Set-TextCell $q4Ws.Cells.Item($r, 7) $row.G   # held value (100M)
    $q4Ws.Cells.Item($r, 8).Value = $row.H        # position rank (numeric)
}

# ---------------------------------------------------------------------
# 2. Insert the new summary row into "总计" (row 2, pushing old data down)
#    and repair the running index column (A) so it stays 0..8 in order.
# ---------------------------------------------------------------------
$totalWs.Rows.Item(2).Insert()

# Copy the (still correctly-styled) row below down into the new blank
# row so column A keeps its index-column style and B:D stay unstyled.
$totalWs.Range("A3:D3").Copy()
$totalWs.Range("A2:D2").PasteSpecial(-4122)

$totalWs.Cells.Item(2, 1).Value = 0
$totalWs.Cells.Item(2, 2).Value = "2022-Q4"
$totalWs.Cells.Item(2, 3).Value = 10
$totalWs.Cells.Item(2, 4).Value = 3.04

# Rows 3..10 used to be rows 2..9 before the insert, so their index
# column (0..7) now needs to shift up to 1..8.
for ($r = 3; $r -le 10; $r++) {
    $totalWs.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# 3. Restore the originally active tab (the workbook was last viewed on
#    "2020-Q4"), since adding/copying sheets changes Excel's selection.
# ---------------------------------------------------------------------
$lastWs.Activate()
